# Generate Report for Handback
# Updates the handback-status workbook with refreshed handoff/handback
# timestamps and priority for the 21a27402-... item, plus the refreshed
# "Latest HO Xliff Generate Date" for both 21a27402 and 3c196807 items on
# the Overview sheet (they share the same generate-date value).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-03 10:18:20"
$wsOverview.Range("G4").Value = "2016-09-03 10:18:20"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
# Priority column (E) for 21a27402 row (3) and 3c196807 row (4)
$wsZh.Range("E3").Value = "mt"
$wsZh.Range("E4").Value = "mt"
# Correspond Handoff Datetime column (H) for 21a27402 row (3) and 3c196807 row (4)
$wsZh.Range("H3").Value = "2016-09-03 10:18:15"
$wsZh.Range("H4").Value = "2016-09-03 10:18:15"
# Correspond Handback DateTime column (K) for 21a27402 row (3) and 3c196807 row (4)
$wsZh.Range("K3").Value = "2016-09-03 10:18:33"
$wsZh.Range("K4").Value = "2016-09-03 10:18:33"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
# Priority column (E) for 21a27402 row (3) and 3c196807 row (4)
$wsDe.Range("E3").Value = "mt"
$wsDe.Range("E4").Value = "mt"
# Correspond Handoff Datetime column (H) for 21a27402 row (3) and 3c196807 row (4)
$wsDe.Range("H3").Value = "2016-09-03 10:18:20"
$wsDe.Range("H4").Value = "2016-09-03 10:18:20"
# Correspond Handback DateTime column (K) for 21a27402 row (3) and 3c196807 row (4)
$wsDe.Range("K3").Value = "2016-09-03 10:18:40"
$wsDe.Range("K4").Value = "2016-09-03 10:18:40"
